$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the effect script bodies (D column) -----------------------------
# Row 2: "room-enter" cue's argument list moved from a positional array
# ([...]) to a keyed dictionary ({"roomname": ...}).
$d2 = "displayCD(""room-enter"", {""roomname"":variableMap[""CURRENT_ROOM""].getRoomName()})`nvariableMap[""JUST_ENTERED""] =0"
$ws.Range("D2").Value = $d2

# Row 3: same change, but the author left an extra space before the brace.
$d3 = "displayCD(""room-enter"",  {""roomname"":variableMap[""CURRENT_ROOM""].getRoomName()})`nvariableMap[""JUST_ENTERED""] =0"
$ws.Range("D3").Value = $d3

# Row 4: "game-start" cue's empty argument list becomes an empty dictionary.
$ws.Range("D4").Value = "displayCD( ""game-start"",{})"

# --- Cosmetic follow-on from the wider text now living in column D ---------
$ws.Columns.Item(4).ColumnWidth = 80.5

# --- Reset the saved selection/view back to the top-left cell --------------
$ws.Range("A1").Select() | Out-Null
